# "ticker appears with results" - append a new results row (row 3) to the
# BIIB sentiment/analysis sheet, carrying over the formats used by row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 42632.882164351853

$ws.Range("B3").Value = 5
$ws.Range("C3").Value = "Neutral"
$ws.Range("D3").Value = 14
$ws.Range("E3").Value = 12000
$ws.Range("F3").Value = 1284
$ws.Range("G3").Value = 58
$ws.Range("H3").Value = 39
$ws.Range("I3").Value = 70
$ws.Range("J3").Value = 29
$ws.Range("K3").Value = 7229
$ws.Range("L3").Value = 151
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 12
$ws.Range("O3").Value = 5
$ws.Range("P3").Value = "Bag"
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 1.77

$ws.Range("S3").Value = 0.1132
$ws.Range("S3").NumberFormat = "0.00%"

$ws.Range("T3").Value = -4.05
$ws.Range("U3").Value = 5.85
$ws.Range("V3").Value = "N/A"
$ws.Range("W3").Value = 0

# New "Neutral" verdict widened column C's best-fit width.
$ws.Columns("C:C").ColumnWidth = 6.83
